# Add product import logic: populate rows 2 and 3 with new product data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: 8000-BW-XLS product
$ws.Range("A2").Value = "8000-BW-XLS"
$ws.Range("B2").Value = 8000
$ws.Range("C2").Value = 11.2
$ws.Range("D2").Value = 12.3
$ws.Range("E2").Value = 110
$ws.Range("F2").Value = 110
$ws.Range("G2").Value = 120
$ws.Range("H2").Value = 15200
$ws.Range("I2").Value = 65.5
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 120
$ws.Range("L2").Value = "9x1231.11"
$ws.Range("M2").Value = "9x1231.11"
$ws.Range("N2").Value = "产品"

# Row 3: 9000-BW-XLS product
$ws.Range("A3").Value = "9000-BW-XLS"
$ws.Range("B3").Value = 9000
$ws.Range("C3").Value = 11.2
$ws.Range("D3").Value = 12.3
$ws.Range("E3").Value = 110
$ws.Range("F3").Value = 110
$ws.Range("G3").Value = 120
$ws.Range("H3").Value = 15200
$ws.Range("I3").Value = 65.5
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 120
$ws.Range("L3").Value = "9x1231.11"
$ws.Range("M3").Value = "9x1231.11"
$ws.Range("N3").Value = "产品"

# Re-apply the data validation list on column N now that rows 2/3 carry values.
$dv = $ws.Range("N1:N1048576")
$dv.Validation.Delete()
$dv.Validation.Add(3, 1, 1, '"产品,配件"')

# Move the active selection to H17, matching the recorded view state.
$ws.Range("H17").Select()

$wb.Save()
